$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# misinterpreted as numbers by Excel (so they stay exact text, matching the
# original inlineStr cell content).
$textForceCells = @("D5", "D7", "D8", "D9", "D10", "D12", "D13", "D15", "D16", "D20", "D21", "D24", "D25", "D26", "D27", "D28", "D29", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D49", "D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cryptocurrency values scraped on Tue Apr 18 21:48:41 UTC 2023
$ws.Range("D2").Value = '30.377.30'
$ws.Range("E2").Value = '  +2.22%  '
$ws.Range("D3").Value = '2.092.98'
$ws.Range("E3").Value = '  -0.10%  '
$ws.Range("E4").Value = '  -0.77%  '
$ws.Range("D5").Value = '342.77'
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("E6").Value = '  -0.71%  '
$ws.Range("D7").Value = '0.5234'
$ws.Range("D8").Value = '0.4422'
$ws.Range("E8").Value = '  +1.03%  '
$ws.Range("D9").Value = '54.45'
$ws.Range("E9").Value = '  +2.76%  '
$ws.Range("D10").Value = '0.09324'
$ws.Range("E10").Value = '  +0.81%  '
$ws.Range("E11").Value = '  +0.43%  '
$ws.Range("D12").Value = '24.77'
$ws.Range("E12").Value = '  -0.51%  '
$ws.Range("D13").Value = '8.609'
$ws.Range("E13").Value = '  +3.55%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '2.110.58'
$ws.Range("E14").Value = '  +0.34%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '6.905'
$ws.Range("E15").Value = '  +2.44%  '
$ws.Range("D16").Value = '101.53'
$ws.Range("E16").Value = '  +2.16%  '
$ws.Range("E17").Value = '  +0.67%  '
$ws.Range("E18").Value = '  -0.70%  '
$ws.Range("E19").Value = '  +1.73%  '
$ws.Range("D20").Value = '0.06665'
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("D21").Value = '6.334'
$ws.Range("E21").Value = '  +2.42%  '
$ws.Range("E22").Value = '  -0.73%  '
$ws.Range("D23").Value = '30.390.80'
$ws.Range("E23").Value = '  +2.14%  '
$ws.Range("D24").Value = '12.53'
$ws.Range("E24").Value = '  +0.28%  '
$ws.Range("D25").Value = '2.314'
$ws.Range("E25").Value = '  -0.21%  '
$ws.Range("D26").Value = '21.82'
$ws.Range("E26").Value = '  -0.55%  '
$ws.Range("D27").Value = '162.99'
$ws.Range("E27").Value = '  +1.02%  '
$ws.Range("D28").Value = '2.504'
$ws.Range("E28").Value = '  -0.50%  '
$ws.Range("D29").Value = '133.18'
$ws.Range("E29").Value = '  +0.14%  '
$ws.Range("E30").Value = '  +0.36%  '
$ws.Range("E31").Value = '  -0.42%  '
$ws.Range("D32").Value = '1.658'
$ws.Range("E32").Value = '  +0.64%  '
$ws.Range("D33").Value = '6.835'
$ws.Range("E33").Value = '  +9.30%  '
$ws.Range("D34").Value = '6.252'
$ws.Range("E34").Value = '  +1.69%  '
$ws.Range("D35").Value = '3.857'
$ws.Range("E35").Value = '  -1.97%  '
$ws.Range("D36").Value = '10.13'
$ws.Range("E36").Value = '  -0.57%  '
$ws.Range("D37").Value = '0.02632'
$ws.Range("E37").Value = '  +2.22%  '
$ws.Range("D38").Value = '0.06832'
$ws.Range("E38").Value = '  +2.06%  '
$ws.Range("D39").Value = '0.6988'
$ws.Range("E39").Value = '  +1.65%  '
$ws.Range("E40").Value = '  +1.15%  '
$ws.Range("D41").Value = '1.338'
$ws.Range("E41").Value = '  +0.60%  '
$ws.Range("D42").Value = '0.2211'
$ws.Range("E42").Value = '  -0.67%  '
$ws.Range("D43").Value = '0.6810'
$ws.Range("E43").Value = '  +1.67%  '
$ws.Range("D44").Value = '14.43'
$ws.Range("E44").Value = '  +1.20%  '
$ws.Range("D45").Value = '2.342'
$ws.Range("E45").Value = '  +1.15%  '
$ws.Range("D46").Value = '1.002'
$ws.Range("E46").Value = '  -0.63%  '
$ws.Range("E47").Value = '  +19.02%  '
$ws.Range("D48").Value = '3.633'
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.00000000355'
$ws.Range("E49").Value = '  -1.27%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").Value = '1.237'
$ws.Range("E50").Value = '  +10.86%  '
